$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the customer / product info values in column B -----------------
$ws.Range("B1").Value = "Bangalore"
$ws.Range("B2").Value = "life"
# Keep this one stored as text (matches original quotePrefix-backed style s=2)
$ws.Range("B3").Value = "'934122"
$ws.Range("B4").Value = "Newmen E120360"

# --- Re-style B4 (new font: Segoe UI, 12pt, teal FF36B9CC) ------------------
$ws.Range("B4").ClearFormats()
$ws.Range("B4").Font.Name = "Segoe UI"
$ws.Range("B4").Font.Size = 12
$ws.Range("B4").Font.Color = 13416758

# Row 4 grows slightly taller to fit the larger font
$ws.Rows(4).RowHeight = 17.25

# --- Final selection ends on B11, matching the saved UI state --------------
$ws.Range("B11").Select()
